$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11.07461541334959
$ws.Range("C2").Value = 13.25362651949753
$ws.Range("D2").Value = 5.426202884781771
$ws.Range("F2").Value = 25.43893045904739
$ws.Range("G2").Value = 30.71617546584923
$ws.Range("H2").Value = 14.89119137541703
$ws.Range("K2").Value = 8.404439517062462
$ws.Range("L2").Value = 11.15275193508069
$ws.Range("M2").Value = 14.16860043662403
$ws.Range("O2").Value = 22.90650877922387
# Row 3
$ws.Range("B3").Value = 10.80827132442156
$ws.Range("C3").Value = 13.29243649415632
$ws.Range("D3").Value = 5.357019768673666
$ws.Range("F3").Value = 25.49642236383876
$ws.Range("G3").Value = 30.81707362864761
$ws.Range("H3").Value = 14.93943122608546
$ws.Range("K3").Value = 8.165488985032766
$ws.Range("L3").Value = 11.16234124483018
$ws.Range("M3").Value = 14.12770537130639
$ws.Range("O3").Value = 22.98836933419144
# Row 4
$ws.Range("B4").Value = 10.64294795579683
$ws.Range("C4").Value = 13.3175746228034
$ws.Range("D4").Value = 5.313455684953329
$ws.Range("F4").Value = 25.53789072787766
$ws.Range("G4").Value = 30.8877946595122
$ws.Range("H4").Value = 14.97120297090932
$ws.Range("K4").Value = 8.013773288976477
$ws.Range("L4").Value = 11.16994992461506
$ws.Range("M4").Value = 14.10456664314422
$ws.Range("O4").Value = 23.04304779728473
# Row 5
$ws.Range("B5").Value = 10.57522566066958
$ws.Range("C5").Value = 13.32814857841056
$ws.Range("D5").Value = 5.295439496132818
$ws.Range("F5").Value = 25.55633747862251
$ws.Range("G5").Value = 30.91881136557175
$ws.Range("H5").Value = 14.98469173806409
$ws.Range("K5").Value = 7.950747919943596
$ws.Range("L5").Value = 11.17348369495393
$ws.Range("M5").Value = 14.09563967978406
$ws.Range("O5").Value = 23.06643896861422
# Row 6
$ws.Range("B6").Value = 10.56396212233794
$ws.Range("C6").Value = 13.32992433108013
$ws.Range("D6").Value = 5.292432289836948
$ws.Range("F6").Value = 25.55949395196035
$ws.Range("G6").Value = 30.92409414217785
$ws.Range("H6").Value = 14.98696425361692
$ws.Range("K6").Value = 7.940211891456022
$ws.Range("L6").Value = 11.17409664884048
$ws.Range("M6").Value = 14.09418788478567
$ws.Range("O6").Value = 23.07039001774442
# Row 7
$ws.Range("B7").Value = 10.64203592097632
$ws.Range("C7").Value = 13.31771588966263
$ws.Range("D7").Value = 5.313213766451383
$ws.Range("F7").Value = 25.53813324340044
$ws.Range("G7").Value = 30.88820407543625
$ws.Range("H7").Value = 14.97138269211355
$ws.Range("K7").Value = 8.01292808658498
$ws.Range("L7").Value = 11.16999582784418
$ws.Range("M7").Value = 14.10444420893284
$ws.Range("O7").Value = 23.04335876897943
# Row 8
$ws.Range("B8").Value = 10.98320810803745
$ws.Range("C8").Value = 13.26673712731286
$ws.Range("D8").Value = 5.402579403635962
$ws.Range("F8").Value = 25.4574718113322
$ws.Range("G8").Value = 30.74914058785155
$ws.Range("H8").Value = 14.90737795944246
$ws.Range("K8").Value = 8.323119686567853
$ws.Range("L8").Value = 11.15570148031948
$ws.Range("M8").Value = 14.15409483588476
$ws.Range("O8").Value = 22.93381714362106
# Row 9
$ws.Range("B9").Value = 11.63400855941341
$ws.Range("C9").Value = 13.17711075633321
$ws.Range("D9").Value = 5.568813086519274
$ws.Range("F9").Value = 25.34835309025898
$ws.Range("G9").Value = 30.54634361004839
$ws.Range("H9").Value = 14.7989280487306
$ws.Range("K9").Value = 8.88947596209141
$ws.Range("L9").Value = 11.1413008563755
$ws.Range("M9").Value = 14.26679035117713
$ws.Range("O9").Value = 22.75409262062557
# Row 10
$ws.Range("B10").Value = 12.0959620437721
$ws.Range("C10").Value = 13.1175106929909
$ws.Range("D10").Value = 5.684937897264408
$ws.Range("F10").Value = 25.29822715341762
$ws.Range("G10").Value = 30.44038108646024
$ws.Range("H10").Value = 14.72963237229945
$ws.Range("K10").Value = 9.277318034607797
$ws.Range("L10").Value = 11.13899448607957
$ws.Range("M10").Value = 14.3584968644774
$ws.Range("O10").Value = 22.64350447781682
# Row 11
$ws.Range("B11").Value = 12.30163927058946
$ws.Range("C11").Value = 13.09174186351488
$ws.Range("D11").Value = 5.736359266224577
$ws.Range("F11").Value = 25.28196724903172
$ws.Range("G11").Value = 30.40159749386255
$ws.Range("H11").Value = 14.70035841167341
$ws.Range("K11").Value = 9.447139247274976
$ws.Range("L11").Value = 11.1397320423058
$ws.Range("M11").Value = 14.40205009933617
$ws.Range("O11").Value = 22.59786757346065
# Row 12
$ws.Range("B12").Value = 12.37880908874359
$ws.Range("C12").Value = 13.08217618098619
$ws.Range("D12").Value = 5.755621053132567
$ws.Range("F12").Value = 25.27675166728151
$ws.Range("G12").Value = 30.38827099116866
$ws.Range("H12").Value = 14.68959626798144
$ws.Range("K12").Value = 9.510460790907429
$ws.Range("L12").Value = 11.14026722261675
$ws.Range("M12").Value = 14.41879738859934
$ws.Range("O12").Value = 22.58125876475995
# Row 13
$ws.Range("B13").Value = 12.36222212175141
$ws.Range("C13").Value = 13.08422777660055
$ws.Range("D13").Value = 5.751482169149915
$ws.Range("F13").Value = 25.27783304580916
$ws.Range("G13").Value = 30.39108052478601
$ws.Range("H13").Value = 14.69189971270191
$ws.Range("K13").Value = 9.496867744265812
$ws.Range("L13").Value = 11.14014059939102
$ws.Range("M13").Value = 14.41517938749829
$ws.Range("O13").Value = 22.5848058209182
# Row 14
$ws.Range("B14").Value = 12.3080028200122
$ws.Range("C14").Value = 13.09095103725402
$ws.Range("D14").Value = 5.737948208091213
$ws.Range("F14").Value = 25.28151928394145
$ws.Range("G14").Value = 30.4004738287727
$ws.Range("H14").Value = 14.69946652555776
$ws.Range("K14").Value = 9.452368700938859
$ws.Range("L14").Value = 11.13977095003628
$ws.Range("M14").Value = 14.40342286036406
$ws.Range("O14").Value = 22.59648766346248
# Row 15
$ws.Range("B15").Value = 12.27469665547148
$ws.Range("C15").Value = 13.09509426338592
$ws.Range("D15").Value = 5.729630626092237
$ws.Range("F15").Value = 25.28389986423308
$ws.Range("G15").Value = 30.40640476346756
$ws.Range("H15").Value = 14.7041435123947
$ws.Range("K15").Value = 9.424982348082017
$ws.Range("L15").Value = 11.13957781966806
$ws.Range("M15").Value = 14.39625452560616
$ws.Range("O15").Value = 22.60373079903984
# Row 16
$ws.Range("B16").Value = 12.08242425754291
$ws.Range("C16").Value = 13.11922171465854
$ws.Range("D16").Value = 5.681548402000652
$ws.Range("F16").Value = 25.29942151684599
$ws.Range("G16").Value = 30.44310571272399
$ws.Range("H16").Value = 14.73159073844418
$ws.Range("K16").Value = 9.266083680545121
$ws.Range("L16").Value = 11.13898215197283
$ws.Range("M16").Value = 14.35568670396525
$ws.Range("O16").Value = 22.64658104756752
# Row 17
$ws.Range("B17").Value = 11.96327311249249
$ws.Range("C17").Value = 13.13436667482748
$ws.Range("D17").Value = 5.65168571265648
$ws.Range("F17").Value = 25.31062005341097
$ws.Range("G17").Value = 30.46803720069446
$ws.Range("H17").Value = 14.74900464665192
$ws.Range("K17").Value = 9.166884967283279
$ws.Range("L17").Value = 11.13907364792282
$ws.Range("M17").Value = 14.33126324826588
$ws.Range("O17").Value = 22.67406539389895
# Row 18
$ws.Range("B18").Value = 11.89432400639739
$ws.Range("C18").Value = 13.14320416179816
$ws.Range("D18").Value = 5.634377776862066
$ws.Range("F18").Value = 25.31767695938586
$ws.Range("G18").Value = 30.48326343263573
$ws.Range("H18").Value = 14.75923235910087
$ws.Range("K18").Value = 9.109208187649127
$ws.Range("L18").Value = 11.13929452010452
$ws.Range("M18").Value = 14.31738886601721
$ws.Range("O18").Value = 22.69031316696202
# Row 19
$ws.Range("B19").Value = 11.87090981357885
$ws.Range("C19").Value = 13.14621813474433
$ws.Range("D19").Value = 5.628495248864704
$ws.Range("F19").Value = 25.32017203661697
$ws.Range("G19").Value = 30.48857081688256
$ws.Range("H19").Value = 14.76273164943248
$ws.Range("K19").Value = 9.089574429928552
$ws.Range("L19").Value = 11.13939823001416
$ws.Range("M19").Value = 14.31272128889408
$ws.Range("O19").Value = 22.69588982788067
# Row 20
$ws.Range("B20").Value = 11.9760006124638
$ws.Range("C20").Value = 13.13274138074131
$ws.Range("D20").Value = 5.654878340192025
$ws.Range("F20").Value = 25.30936420876453
$ws.Range("G20").Value = 30.46529143045668
$ws.Range("H20").Value = 14.74712899793676
$ws.Range("K20").Value = 9.177509308409162
$ws.Range("L20").Value = 11.13904650238759
$ws.Range("M20").Value = 14.33384529285336
$ws.Range("O20").Value = 22.67109414154583
# Row 21
$ws.Range("B21").Value = 12.32394832129806
$ws.Range("C21").Value = 13.08897103677224
$ws.Range("D21").Value = 5.741929235203745
$ws.Range("F21").Value = 25.28041098484564
$ws.Range("G21").Value = 30.39767783697376
$ws.Range("H21").Value = 14.69723519676704
$ws.Range("K21").Value = 9.465466174392896
$ws.Range("L21").Value = 11.13987258932937
$ws.Range("M21").Value = 14.40686920278403
$ws.Range("O21").Value = 22.59303815053459
# Row 22
$ws.Range("B22").Value = 12.54715084689246
$ws.Range("C22").Value = 13.06148576554028
$ws.Range("D22").Value = 5.797592183262941
$ws.Range("F22").Value = 25.26697731381272
$ws.Range("G22").Value = 30.3614176581992
$ws.Range("H22").Value = 14.66651092119108
$ws.Range("K22").Value = 9.647901258313825
$ws.Range("L22").Value = 11.14190340197616
$ws.Range("M22").Value = 14.45607415049596
$ws.Range("O22").Value = 22.54594679089945
# Row 23
$ws.Range("B23").Value = 12.42843057477132
$ws.Range("C23").Value = 13.0760528408205
$ws.Range("D23").Value = 5.767999039205208
$ws.Range("F23").Value = 25.27364472539244
$ws.Range("G23").Value = 30.38004322602563
$ws.Range("H23").Value = 14.68273668203964
$ws.Range("K23").Value = 9.551070150869657
$ws.Range("L23").Value = 11.14068347648138
$ws.Range("M23").Value = 14.42968032052635
$ws.Range("O23").Value = 22.57072097646897
# Row 24
$ws.Range("B24").Value = 11.97024789804527
$ws.Range("C24").Value = 13.13347576980383
$ws.Range("D24").Value = 5.653435387083569
$ws.Range("F24").Value = 25.3099300489457
$ws.Range("G24").Value = 30.46653001264035
$ws.Range("H24").Value = 14.7479763051813
$ws.Range("K24").Value = 9.172708053585655
$ws.Range("L24").Value = 11.13905825066646
$ws.Range("M24").Value = 14.33267742964518
$ws.Range("O24").Value = 22.6724360535273
# Row 25
$ws.Range("B25").Value = 11.46044714370389
$ws.Range("C25").Value = 13.20025570784327
$ws.Range("D25").Value = 5.524861437751032
$ws.Range("F25").Value = 25.3726047093644
$ws.Range("G25").Value = 30.5936764981458
$ws.Range("H25").Value = 14.82644195961541
$ws.Range("K25").Value = 8.741046554239825
$ws.Range("L25").Value = 11.14374061543861
$ws.Range("M25").Value = 14.23470549026356
$ws.Range("O25").Value = 22.79894971704186

Write-Host "Updated loading_percent values for 380 kV case"
